# Generate Report for Handback
# ------------------------------------------------------------------
# This script reproduces, via Excel COM automation, the "handback"
# report-generation edit: the two target files (identified by their
# GUID-prefixed names) have now been handed back from localization,
# so each per-locale sheet (zh-cn, de-de) gets its "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns
# populated, the Overview sheet's summary status flips from
# "Ready for handoff" to "Handed back: in sync with en-US", and a
# few columns are widened so the new content is readable.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7eb7d22c30c52c8d2a97bdf0cf4e527094dba8af/e2e/"
$file1 = "2a357830-f0a0-4343-9c6d-5932f88482b2.md"
$file2 = "37be8748-e7d1-4d61-9d91-3f01f47bc0af.md"

# ------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by the Overview sheet (zh-cn/de-de summary
#    columns E/F) and by the per-locale sheets' "Status" column (C).
#    Update every cell that carries it so the underlying shared
#    string is fully replaced (no leftover "Ready for handoff" cell).
# ------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ------------------------------------------------------------------
# 2. Populate "Latest Target File" (col I) / "Latest Handback File"
#    (col J) / "Latest Handback DateTime" (col K) for both rows on
#    both locale sheets.
# ------------------------------------------------------------------

# --- zh-cn ---
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($githubBase + $file1), "", "", $file1) | Out-Null
$wsZh.Range("J2").Value = "2a357830-f0a0-4343-9c6d-5932f88482b2.a60ff4a15651c53f48bb04d47ed0d98d2763d353.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-02 22:33:54"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($githubBase + $file2), "", "", $file2) | Out-Null
$wsZh.Range("J3").Value = "37be8748-e7d1-4d61-9d91-3f01f47bc0af.496cba0b6da7a09fa56ad63c7185fd969eb6874f.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-02 22:33:54"

# --- de-de ---
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($githubBase + $file1), "", "", $file1) | Out-Null
$wsDe.Range("J2").Value = "2a357830-f0a0-4343-9c6d-5932f88482b2.a60ff4a15651c53f48bb04d47ed0d98d2763d353.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-02 22:34:03"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($githubBase + $file2), "", "", $file2) | Out-Null
$wsDe.Range("J3").Value = "37be8748-e7d1-4d61-9d91-3f01f47bc0af.496cba0b6da7a09fa56ad63c7185fd969eb6874f.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-02 22:34:03"

# ------------------------------------------------------------------
# 3. Widen columns so the newly-filled-in long file names / links are
#    readable: Overview E/F, and col C (Status) + I/J (Target /
#    Handback file) on both locale sheets.
# ------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

$wsZh.Columns.Item(3).ColumnWidth = 29.17
$wsZh.Columns.Item(9).ColumnWidth = 39.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

$wsDe.Columns.Item(3).ColumnWidth = 29.17
$wsDe.Columns.Item(9).ColumnWidth = 39.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1
